# Fixed logger merge issue in NormalParser and Ui
#
# 1. Slide master background: blue (0000FF) -> navy (000080)
# 2. All white (FFFFFF) title/body text runs on the content slides -> sienna (A0522D)

$p = $ppt.ActivePresentation

# --- Slide master background color ---
$master = $p.SlideMaster
$master.Background.Fill.ForeColor.RGB = 8388608   # RGB(0, 0, 128) = 0x000080

# --- Text run colors across all slides ---
$oldColor = 16777215   # RGB(255, 255, 255) = 0xFFFFFF
$newColor = 2970272    # RGB(160, 82, 45)   = 0xA0522D

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $textRange = $shape.TextFrame.TextRange
            if ($textRange.Font.Color.RGB -eq $oldColor) {
                $textRange.Font.Color.RGB = $newColor
            }
        }
    }
}
